# Update mods data [2025-12-16 15:13:09]
# Append a new row (37) to the ModCounts sheet with the latest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (36) down onto
# the new row 37 so the new cells keep the table's center/middle alignment.
$ws.Range("A36:C36").Copy()
$ws.Range("A37:C37").PasteSpecial(-4122)  # xlPasteFormats

# Write the new row's values. The date column stores the date as literal
# text (matching the rest of the column, which is plain text, not a real
# Excel date) - prefix with an apostrophe so it is entered as text instead
# of being auto-converted into a date serial number.
$ws.Range("A37").Value = "'2025/12/16"
$ws.Range("B37").Value = "逃离鸭科夫"
$ws.Range("C37").Value = 1343
